$d = $word.ActiveDocument

# Replace the merge-field property reference "externalShortName" with the
# snake_case "external_short_name" so the template now reads
# <<caseManagementLocation.external_short_name>> instead of
# <<caseManagementLocation.externalShortName>>.
$d.Content.Find.Execute("externalShortName", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "external_short_name", 2) | Out-Null
